# "Testing if push works" - update the Codebook sheet's allowed-values
# text for HC_time (row 5) from the old "-999" sentinel wording to the
# "NA" wording already used elsewhere in the column, and move the
# sheet's active selection from C6 to G9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebook")

# Make sure we're working on/viewing the right sheet, matching
# tabSelected="1" on this sheet in the saved file.
$ws.Activate()

# C5: "numeric value >0 or -999" -> "numeric value >0 or NA"
$ws.Range("C5").Value = "numeric value >0 or NA"

# Move the selection shown when the sheet is reopened.
$ws.Range("G9").Select()
